$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.690722107887268
$ws.Range("B1").Value = 5.390152931213379
$ws.Range("C1").Value = 2.733572959899902
$ws.Range("D1").Value = 2.370488405227661
$ws.Range("E1").Value = 2.159883260726929
